$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Jurisdiction-Basic" form section contained an accidental duplicate
# "Name of Jurisdiction" / STRING field directly above each "Jurisdiction" /
# JURISDICTION field (rows 29-30 and rows 32-33 of the original sheet).
# Remove the two duplicate/erroneous rows so only the real JURISDICTION
# field-type row remains for each occurrence.
$ws.Rows("29").Delete()
$ws.Rows("31").Delete()

# The last leftover "Name of Jurisdiction" / STRING row (originally row 38,
# now row 35 after the two deletions above) was a stray, non-functional
# field definition left at the end of the form. Turn it into a proper,
# fully-functional "Jurisdiction" / JURISDICTION field, matching the other
# Jurisdiction entries above.
$ws.Range("C35").Value = "Jurisdiction"
$ws.Range("D35").Value = "JURISDICTION"

# Reflect where editing ended up, matching the author's final selection.
$ws.Range("C32").Select()
